# ResonatorArray.xlsx fix-up:
#  - corrected "Wafer Shift" x/y coordinates (cols M,N,O,P) on a number of
#    mask/cell rows (array positioning fix)
#  - extended the sheet up to row 1 (adds blank rows 1-10 above the data
#    that used to start at row 11)
#  - added explicit column widths for columns A:Y

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Extend the sheet's used range up to A1 and materialize blank rows
#    1-10 above the existing data (which starts at row 11), without
#    disturbing any of the existing rows/data.
# ---------------------------------------------------------------------

# Touch column A of the new top rows so the workbook recomputes its
# used range/dimension to start at A1, and materializes row stubs for
# rows 1-10 (copy formatting only, from an already-styled cell, so no
# new style entries are introduced and no cell values are set).
$ws.Range("A11").Copy()
for ($r = 1; $r -le 10; $r++) {
    $ws.Range("A" + $r).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 2) Apply explicit column widths (A:Y)
# ---------------------------------------------------------------------
$colWidths = @(16, 6, 29, 7, 8, 6, 6, 21, 21, 22, 21, 6, 19, 9, 20, 20, 6, 6, 24, 6, 6, 9, 8, 7, 8)
for ($c = 1; $c -le $colWidths.Length; $c++) {
    # Excel's ColumnWidth property is expressed in "characters" and is
    # internally offset by 5/6 of a character from the stored OOXML
    # width, so back that offset out to land exactly on the desired
    # stored width.
    $ws.Columns.Item($c).ColumnWidth = ($colWidths[$c - 1] - 5 / 6)
}

# ---------------------------------------------------------------------
# 3) Corrected Wafer Shift coordinates (columns M=x, N=y, O=x, P=y)
# ---------------------------------------------------------------------
$ws.Range("N15").Value = -38.401
$ws.Range("P15").Value = 38.251

$ws.Range("N16").Value = -37.442
$ws.Range("P16").Value = 48.842

$ws.Range("M17").Value = 5.992000000000001
$ws.Range("N17").Value = -16.296
$ws.Range("O17").Value = -2.992000000000001
$ws.Range("P17").Value = 16.946

$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -38.5
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 48.7

$ws.Range("M19").Value = 5.992
$ws.Range("N19").Value = -27.496
$ws.Range("O19").Value = -16.592
$ws.Range("P19").Value = 27.496

$ws.Range("N20").Value = -38.401
$ws.Range("P20").Value = 38.401

$ws.Range("M21").Value = -88
$ws.Range("N21").Value = -31.082
$ws.Range("O21").Value = 79.2
$ws.Range("P21").Value = 31.082

$ws.Range("M22").Value = 88
$ws.Range("N22").Value = -20.082
$ws.Range("O22").Value = -96.8
$ws.Range("P22").Value = 20.082

$ws.Range("M24").Value = -99.90000000000001
$ws.Range("O24").Value = 108.55

$ws.Range("N25").Value = -39.168
$ws.Range("P25").Value = 39.568

$ws.Range("N26").Value = -38.401
$ws.Range("P26").Value = 36.251

$ws.Range("N27").Value = -50.55
$ws.Range("P27").Value = 63.15

$ws.Range("M28").Value = -50.55
$ws.Range("O28").Value = 40.55

$ws.Range("N93").Value = -38.401
$ws.Range("P93").Value = 38.401

$ws.Range("N94").Value = -38.401
$ws.Range("P94").Value = 38.801

$ws.Range("N95").Value = -38.401
$ws.Range("P95").Value = 38.401

$ws.Range("M97").Value = -99.90000000000001
$ws.Range("O97").Value = 109.65

$ws.Range("M98").Value = 5.992000000000001
$ws.Range("N98").Value = -16.296
$ws.Range("O98").Value = -2.992000000000001
$ws.Range("P98").Value = 15.946

$ws.Range("N99").Value = -37.645
$ws.Range("P99").Value = 38.295

$ws.Range("N100").Value = -39.063
$ws.Range("P100").Value = 38.813

$ws.Range("N101").Value = -37.442
$ws.Range("P101").Value = 49.442

$ws.Range("M102").Value = 0
$ws.Range("N102").Value = -38.5
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 49.3

$ws.Range("M103").Value = 5.992
$ws.Range("N103").Value = -27.496
$ws.Range("O103").Value = 4.608
$ws.Range("P103").Value = 27.496

$ws.Range("M104").Value = -88
$ws.Range("N104").Value = -31.082
$ws.Range("O104").Value = 78.59999999999999
$ws.Range("P104").Value = 31.082

$ws.Range("M105").Value = 88
$ws.Range("N105").Value = -20.082
$ws.Range("O105").Value = -97.40000000000001
$ws.Range("P105").Value = 20.082
